# Weekly update: insert a new week's price record for Membrillo (Champion,
# Primera) above the current row 51, shifting the existing rows 51-66 down
# to 52-67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 51 (pushes rows 51:66 down to 52:67).
$ws.Rows.Item(51).Insert()

# Fill the new row 51 with the new week's record.
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 45127
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100104
$ws.Cells.Item(51, 8).Value = "Frutos de pepita"
$ws.Cells.Item(51, 9).Value = 100104003
$ws.Cells.Item(51, 10).Value = "Membrillo"
$ws.Cells.Item(51, 11).Value = "Champion"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 80
$ws.Cells.Item(51, 14).Value = 12000
$ws.Cells.Item(51, 15).Value = 12000
$ws.Cells.Item(51, 16).Value = 12000
$ws.Cells.Item(51, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 667
$ws.Cells.Item(51, 20).Value = 18
